# Update gh-pages to output generated at 456a3b4
# Applies refreshed "想去人数" (want-to-go count) and "最低票价" (min ticket
# price) figures to the 展览 (Exhibitions) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("G3").Value = 38
$wsExhibitions.Range("G4").Value = 30
$wsExhibitions.Range("F9").Value = 1187
$wsExhibitions.Range("F10").Value = 278
$wsExhibitions.Range("F11").Value = 337
$wsExhibitions.Range("F12").Value = 10365
$wsExhibitions.Range("F14").Value = 78
$wsExhibitions.Range("F16").Value = 342
$wsExhibitions.Range("F17").Value = 660
$wsExhibitions.Range("F18").Value = 11904
$wsExhibitions.Range("F19").Value = 12295

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("G3").Value = 38
$wsAllTypes.Range("G4").Value = 30
$wsAllTypes.Range("F10").Value = 1187
$wsAllTypes.Range("F12").Value = 337
$wsAllTypes.Range("F13").Value = 10365
$wsAllTypes.Range("F15").Value = 78
$wsAllTypes.Range("F17").Value = 343
$wsAllTypes.Range("F18").Value = 660
$wsAllTypes.Range("F19").Value = 11904
$wsAllTypes.Range("F20").Value = 12295
